$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Marker Student ..." headers to the new "...(Marker Student)" form.
$ws.Range("F1").Value = "Matriculation Number(Marker Student)"
$ws.Range("G1").Value = "First Name(Marker Student)"
$ws.Range("H1").Value = "Last Name(Marker Student)"
$ws.Range("I1").Value = "Email(Marker Student)"

# Re-fit the affected columns so their widths track the new (slightly longer)
# header text, as Excel's bestFit auto-sizing would do after the rename.
$ws.Columns.Item(6).ColumnWidth = 35.736979166666664
$ws.Columns.Item(7).ColumnWidth = 25.451822916666668
$ws.Columns.Item(8).ColumnWidth = 25.022135416666668
$ws.Columns.Item(9).ColumnWidth = 20.592447916666668
